$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handback
#
# 1) The "Status" text ("Ready for handoff") becomes "Handed back: in sync
#    with en-US" everywhere it is used (Overview!B2:C3, zh-cn!C2:C3,
#    de-de!C2:C3 all share the same underlying string).
# 2) On the zh-cn and de-de detail sheets, two new columns are populated for
#    each data row:
#       F = "Latest Target File"   -> same file/link as column A (source .md)
#       G = "Latest Handback File" -> same file/link as column D (target .xlf)
# 3) The "Latest Handback DateTime" column (H) is filled in with the actual
#    handback timestamp (it previously held the zero-date placeholder).
# ---------------------------------------------------------------------------

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$newStatus = "Handed back: in sync with en-US"

# --- 1) Status column everywhere -------------------------------------------
$ws1.Range("B2").Value = $newStatus
$ws1.Range("C2").Value = $newStatus
$ws1.Range("B3").Value = $newStatus
$ws1.Range("C3").Value = $newStatus

$ws2.Range("C2").Value = $newStatus
$ws2.Range("C3").Value = $newStatus

$ws3.Range("C2").Value = $newStatus
$ws3.Range("C3").Value = $newStatus

# Cornflower-blue underlined hyperlink-style font (matches the workbook's
# existing custom "HyperLink" cell style used in columns A/B/D).
$hyperlinkColor = 15570276   # RGB(0x64,0x95,0xED) packed as BGR OLE color
$hyperlinkUnderline = 2      # xlUnderlineStyleSingle

function Set-HandbackLinkCell {
    param($ws, $cellRef, $displayText, $address)

    $range = $ws.Range($cellRef)
    $range.Font.Underline = $hyperlinkUnderline
    $range.Font.Color = $hyperlinkColor
    $ws.Hyperlinks.Add($range, $address, "", "", $displayText) | Out-Null
}

# --- 2) zh-cn sheet: fill F (Latest Target File) / G (Latest Handback File) -
Set-HandbackLinkCell $ws2 "F2" "68675ab0-7c39-48f7-bb88-c13f0a2273f4.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/675c75d6a198c176b6ed29ad97c5bd80980ee1d7/e2e/68675ab0-7c39-48f7-bb88-c13f0a2273f4.md"
Set-HandbackLinkCell $ws2 "G2" "68675ab0-7c39-48f7-bb88-c13f0a2273f4.de8e460e3130a060c6aa4806d7a7b79258a26fad.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cf56087067d144a2f4eb656d7eb1b53cbd4d94d7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/68675ab0-7c39-48f7-bb88-c13f0a2273f4.de8e460e3130a060c6aa4806d7a7b79258a26fad.zh-cn.xlf"

Set-HandbackLinkCell $ws2 "F3" "a3bad4c8-8470-4fbc-813c-39fb757c7186.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/675c75d6a198c176b6ed29ad97c5bd80980ee1d7/e2e/a3bad4c8-8470-4fbc-813c-39fb757c7186.md"
Set-HandbackLinkCell $ws2 "G3" "a3bad4c8-8470-4fbc-813c-39fb757c7186.17bb28e19de1fa4954dd62f3c1359502986c9a4e.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cf56087067d144a2f4eb656d7eb1b53cbd4d94d7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a3bad4c8-8470-4fbc-813c-39fb757c7186.17bb28e19de1fa4954dd62f3c1359502986c9a4e.zh-cn.xlf"

# Latest Handback DateTime (H) now carries the real handback timestamp.
$ws2.Range("H2").Value = "2016-03-21 06:54:55"
$ws2.Range("H3").Value = "2016-03-21 06:54:55"

# --- 3) de-de sheet: fill F (Latest Target File) / G (Latest Handback File) -
Set-HandbackLinkCell $ws3 "F2" "68675ab0-7c39-48f7-bb88-c13f0a2273f4.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/675c75d6a198c176b6ed29ad97c5bd80980ee1d7/e2e/68675ab0-7c39-48f7-bb88-c13f0a2273f4.md"
Set-HandbackLinkCell $ws3 "G2" "68675ab0-7c39-48f7-bb88-c13f0a2273f4.de8e460e3130a060c6aa4806d7a7b79258a26fad.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/012968918c443838d7a4735655941654d11ff7ae/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/68675ab0-7c39-48f7-bb88-c13f0a2273f4.de8e460e3130a060c6aa4806d7a7b79258a26fad.de-de.xlf"

Set-HandbackLinkCell $ws3 "F3" "a3bad4c8-8470-4fbc-813c-39fb757c7186.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/675c75d6a198c176b6ed29ad97c5bd80980ee1d7/e2e/a3bad4c8-8470-4fbc-813c-39fb757c7186.md"
Set-HandbackLinkCell $ws3 "G3" "a3bad4c8-8470-4fbc-813c-39fb757c7186.17bb28e19de1fa4954dd62f3c1359502986c9a4e.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/012968918c443838d7a4735655941654d11ff7ae/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a3bad4c8-8470-4fbc-813c-39fb757c7186.17bb28e19de1fa4954dd62f3c1359502986c9a4e.de-de.xlf"

# Latest Handback DateTime (H) now carries the real handback timestamp.
$ws3.Range("H2").Value = "2016-03-21 06:55:01"
$ws3.Range("H3").Value = "2016-03-21 06:55:01"

Write-Host "Handback report generated."
